$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-78) holds a date serial that must be incremented by one
# day (45204 -> 45205, i.e. 2023-10-05 -> 2023-10-06).
for ($row = 2; $row -le 78; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
